# Quotation_Tool.xlsx edit:
#   - Add a new "BAS_USD" column (J) to the Freight sheet holding the base
#     ocean-freight rate per lane, and expose it as the FR_BAS named range.
#   - Rewrite the Freight/Container cost formulas (Quote!D17, Quote!D18) to
#     use only FR_BAS * 1.0605 instead of (FR_USD + FR_EUR*$B$12) * 1.0605.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New defined name FR_BAS -> Freight!$J$2:$J$141
# ---------------------------------------------------------------------------
$wb.Names.Add('FR_BAS', '=Freight!$J$2:$J$141')

# ---------------------------------------------------------------------------
# 2. Freight sheet: header + 140 BAS_USD values (rows 2-141)
# ---------------------------------------------------------------------------
$frWs = $wb.Worksheets.Item("Freight")

$frWs.Cells.Item(1, 10).Value = "BAS_USD"

$basValues = @(
    2400,2122,1140,790,1290,1450,3250,890,1450,830,
    790,1490,1340,2425,1519,100,1530,1500,1450,1290,
    1400,1750,1960,670,2350,790,1490,1290,1550,1190,
    100,200,2693,100,790,100,1820,790,1510,990,
    1850,790,1450,1290,2350,790,1400,2850,890,348,
    2250,790,1770,790,1490,1090,1,1540,2250,100,
    790,3050,1953,2700,3790,790,1950,790,1450,100,
    2250,1972,990,640,1140,1300,3100,740,1300,550,
    640,1340,1190,2275,1369,12,1380,1350,1300,1140,
    1250,1600,1810,1020,2200,640,1340,1140,1400,1040,
    209,17,2543,353,640,21,1670,640,1360,840,
    1700,640,1300,1140,2200,640,1250,2700,740,571,
    2100,640,2120,640,1340,940,420,1390,2100,168,
    640,2900,1803,2550,3640,640,1800,640,1300,100
)

$row = 2
foreach ($v in $basValues) {
    $frWs.Cells.Item($row, 10).Value = $v
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3. Quote sheet: Freight/Container formulas now reference FR_BAS only
# ---------------------------------------------------------------------------
$quoteWs = $wb.Worksheets.Item("Quote")

$quoteWs.Range("D17").Formula = '=IF(NOT($K$6),"",IF(K9=0,"-",INDEX(FR_BAS,K9)*1.0605))'
$quoteWs.Range("D18").Formula = '=IF(NOT($K$6),"",IF(K10=0,"-",INDEX(FR_BAS,K10)*1.0605))'
